$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Helper: write a text value to a cell without Excel's autodetection turning
# it into a date/number (e.g. "2012-04-30" would otherwise become a date
# serial). We briefly force a text format, assign, then restore the cell
# back to the default "Normal" style so no stray formatting is left behind.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Header row (row 1): fix B1:G1 and extend with H1:N1 ---
$ws.Cells.Item(1, 2).Value = "species"
$ws.Cells.Item(1, 3).Value = "debtor"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "total"
$ws.Cells.Item(1, 6).Value = "register_date"
$ws.Cells.Item(1, 7).Value = "register_reason"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Apply the existing header style (bold, bordered, centered -- same as B1:G1)
# to the new header cells H1:N1 by copying formats from an existing header
# cell. This reuses the workbook's existing style index instead of building
# up new cellXfs entries one property at a time.
$ws.Cells.Item(1, 2).Copy()
$ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(1, 14)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2 (A2 = 85): 現金 / 高金素梅 / 陳麗卿... / 6000000 / ... ---
$ws.Cells.Item(2, 2).Value = "現金"
$ws.Cells.Item(2, 3).Value = "高金素梅"
$ws.Cells.Item(2, 4).Value = "陳麗卿新北市泰山區明志路"
$ws.Cells.Item(2, 5).Value = 6000000
$ws.Cells.Item(2, 6).Value = "96年02月06日"
$ws.Cells.Item(2, 7).Value = "借款"
$ws.Cells.Item(2, 8).Value = "debt"
$ws.Cells.Item(2, 9).Value = "normal"
Set-TextValue $ws.Cells.Item(2, 10) "2012-04-30"
$ws.Cells.Item(2, 11).Value = "高金素梅"
$ws.Cells.Item(2, 12).Value = 926
$ws.Cells.Item(2, 13).Value = "tmp92521"
$ws.Cells.Item(2, 14).Value = 85

# --- Row 3 (A3 = 86): 現金 / 局金素梅 / 石旭松... / 4000000 / ... ---
$ws.Cells.Item(3, 2).Value = "現金"
$ws.Cells.Item(3, 3).Value = "局金素梅"
$ws.Cells.Item(3, 4).Value = "石旭松新北市泰山區明志路"
$ws.Cells.Item(3, 5).Value = 4000000
$ws.Cells.Item(3, 6).Value = "96年02月06日"
$ws.Cells.Item(3, 7).Value = "借款"
$ws.Cells.Item(3, 8).Value = "debt"
$ws.Cells.Item(3, 9).Value = "normal"
Set-TextValue $ws.Cells.Item(3, 10) "2012-04-30"
$ws.Cells.Item(3, 11).Value = "高金素梅"
$ws.Cells.Item(3, 12).Value = 926
$ws.Cells.Item(3, 13).Value = "tmp92521"
$ws.Cells.Item(3, 14).Value = 86
